$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, so the existing NIF/numbers
# column (C) shifts to column D, freeing up column C for the new
# "Email" column.
$ws.Columns("C").Insert()

# New column header
$ws.Range("C1").Value = "Email"

# Email values + hyperlinks (mailto:) for each person
$ws.Range("C2").Value = "juan@gmail.com"
$ws.Range("C3").Value = "luis@gmail.com"
$ws.Range("C4").Value = "ana@gmail.com"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:juan@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:luis@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:ana@gmail.com")

$ws.Range("B6").Select()
